$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing it to remain plain text, so that
# number-looking strings (e.g. "0.998") are not auto-converted to numerics.
# NumberFormat is restored to the default ("Normal" style) afterwards so no
# visible formatting change is introduced.
function Set-TextValue {
    param($addr, $value)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# --- Update Price (D) and Volume(1h) (E) values reflecting the latest scrape ---
$ws.Range("D2").Value = "54.340.08"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "2.286.70"
$ws.Range("E3").Value = "  +2.77%  "
Set-TextValue "D4" "0.998"
$ws.Range("E4").Value = "  -0.24%  "
Set-TextValue "D5" "496.01"
$ws.Range("E5").Value = "  +2.60%  "
Set-TextValue "D6" "127.93"
$ws.Range("E6").Value = "  +2.16%  "
Set-TextValue "D7" "0.998"
$ws.Range("E7").Value = "  -0.15%  "
Set-TextValue "D8" "0.529"
$ws.Range("E8").Value = "  +2.51%  "
$ws.Range("D9").Value = "2.284.40"
$ws.Range("E9").Value = "  +2.12%  "
Set-TextValue "D10" "0.0950"
$ws.Range("E10").Value = "  +4.01%  "
Set-TextValue "D11" "0.152"
$ws.Range("E11").Value = "  +2.48%  "
$ws.Range("E12").Value = "  +4.35%  "
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").Value = "2.689.07"
$ws.Range("E14").Value = "  +2.52%  "
Set-TextValue "D15" "21.82"
$ws.Range("E15").Value = "  +4.10%  "
$ws.Range("D16").Value = "54.256.45"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("D18").Value = "2.313.79"
$ws.Range("E18").Value = "  +3.66%  "
Set-TextValue "D19" "10.05"
$ws.Range("E19").Value = "  +5.37%  "
$ws.Range("E20").Value = "  +3.73%  "
Set-TextValue "D21" "301.18"
$ws.Range("E21").Value = "  +1.16%  "
Set-TextValue "D22" "6.44"
$ws.Range("E22").Value = "  +5.79%  "
Set-TextValue "D23" "0.999"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  -2.57%  "
Set-TextValue "D25" "62.65"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("E26").Value = "  +1.39%  "
Set-TextValue "D27" "0.373"
$ws.Range("E27").Value = "  +3.07%  "
$ws.Range("D28").Value = "2.380.89"
$ws.Range("E28").Value = "  +2.09%  "
$ws.Range("E29").Value = "  +4.14%  "
Set-TextValue "D30" "7.07"
$ws.Range("E30").Value = "  +1.63%  "
Set-TextValue "D31" "168.94"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("E33").Value = "  +1.93%  "
Set-TextValue "D34" "5.86"
$ws.Range("E34").Value = "  +2.80%  "
Set-TextValue "D36" "0.999"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  +2.69%  "
Set-TextValue "D38" "17.70"
$ws.Range("E38").Value = "  +1.83%  "
Set-TextValue "D39" "1.20"
$ws.Range("E39").Value = "  +4.58%  "
Set-TextValue "D40" "0.869"
$ws.Range("E40").Value = "  +5.05%  "
Set-TextValue "D41" "3.71"
$ws.Range("E41").Value = "  +4.20%  "
Set-TextValue "D42" "35.44"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("E43").Value = "  +3.51%  "
Set-TextValue "D44" "0.373"
$ws.Range("E44").Value = "  +2.74%  "
$ws.Range("E45").Value = "  +3.10%  "
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("E49").Value = "  +1.63%  "
Set-TextValue "D50" "238.60"
$ws.Range("E50").Value = "  +4.45%  "
$ws.Range("E51").Value = "  +3.21%  "

# --- Rows 46/47: Aave and RenderToken swapped rank position, with updated values ---
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "127.38"
$ws.Range("E46").Value = "  +4.01%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D47" "4.84"
$ws.Range("E47").Value = "  +5.07%  "
